# Apply the edit described by the diff:
# 1) Swap match data between row 72 and row 73 (columns F,G,H,I,J,L,M,N,P,Q,R,T,U,V)
# 2) Append a new row 75 with the Luzern vs Grasshoppers match
# 3) The sheet dimension will grow to A1:V75 automatically as new cells are written

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Step 1: swap the content of row 72 and row 73
# ---------------------------------------------------------------------------
$cols = @(6,7,8,9,10,12,13,14,16,17,18,20,21,22)  # F,G,H,I,J,L,M,N,P,Q,R,T,U,V

foreach ($c in $cols) {
    $v72 = $ws.Cells.Item(72, $c).Value2
    $v73 = $ws.Cells.Item(73, $c).Value2
    $ws.Cells.Item(72, $c).Value = $v73
    $ws.Cells.Item(73, $c).Value = $v72
}

# ---------------------------------------------------------------------------
# Step 2: append the new row (row 75) for Luzern vs Grasshoppers
# ---------------------------------------------------------------------------
$newRow = 75

# Copy number/cell formatting from the row above so the new row matches the
# existing look (bold/centered/bordered index column, date-formatted E column)
$ws.Cells.Item(74, 1).Copy() | Out-Null
$ws.Cells.Item($newRow, 1).PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$ws.Cells.Item(74, 5).Copy() | Out-Null
$ws.Cells.Item($newRow, 5).PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$ws.Cells.Item($newRow, 1).Value = 74

$ws.Cells.Item($newRow, 2).Value = "switzerland"
$ws.Cells.Item($newRow, 3).Value = "super-league"
$ws.Cells.Item($newRow, 4).Value = "2023-2024"

$ws.Cells.Item($newRow, 5).Value = 45235.59375

$ws.Cells.Item($newRow, 6).Value = "Luzern"
$ws.Cells.Item($newRow, 7).Value = 2
$ws.Cells.Item($newRow, 8).Value = "Grasshoppers"
$ws.Cells.Item($newRow, 9).Value = 0
$ws.Cells.Item($newRow, 10).Value = 1.56
$ws.Cells.Item($newRow, 11).Value = "29/10/2023 16:42"
$ws.Cells.Item($newRow, 12).Value = 1.91
$ws.Cells.Item($newRow, 13).Value = "05/11/2023 14:06"
$ws.Cells.Item($newRow, 14).Value = 4.64
$ws.Cells.Item($newRow, 15).Value = "29/10/2023 16:42"
$ws.Cells.Item($newRow, 16).Value = 3.96
$ws.Cells.Item($newRow, 17).Value = "05/11/2023 14:06"
$ws.Cells.Item($newRow, 18).Value = 5.2
$ws.Cells.Item($newRow, 19).Value = "29/10/2023 16:42"
$ws.Cells.Item($newRow, 20).Value = 3.87
$ws.Cells.Item($newRow, 21).Value = "05/11/2023 14:06"
$ws.Cells.Item($newRow, 22).Value = "https://www.betexplorer.com/football/switzerland/super-league/luzern-grasshoppers/Oby70Cup/"
